# Auto-generated edit script: updates Shinryu_Profits leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW) with refreshed market-price-derived values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 15384645
$ws.Range("I8").Value = 15384645
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 46153935
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("M8").Value = -46153796
$ws.Range("H19").Value = 8431753
$ws.Range("I19").Value = 6708594.5
$ws.Range("J19").Value = 11112222
$ws.Range("K19").Value = 6708594.5
$ws.Range("L19").Value = 11112222
$ws.Range("M19").Value = -6708419.5
$ws.Range("N19").Value = -11112572
$ws.Range("H64").Value = 3553.311
$ws.Range("I64").Value = 3329.6296
$ws.Range("J64").Value = 3888.8333
$ws.Range("K64").Value = 3329.6296
$ws.Range("L64").Value = 3888.8333
$ws.Range("M64").Value = -3081.6296
$ws.Range("N64").Value = -4384.8333
$ws.Range("H67").Value = 3553.311
$ws.Range("I67").Value = 3329.6296
$ws.Range("J67").Value = 3888.8333
$ws.Range("K67").Value = 3329.6296
$ws.Range("L67").Value = 3888.8333
$ws.Range("M67").Value = -2471.6296
$ws.Range("N67").Value = -5604.8333
$ws.Range("H94").Value = 3471.7273
$ws.Range("I94").Value = 3111.25
$ws.Range("J94").Value = 4433
$ws.Range("K94").Value = 3111.25
$ws.Range("L94").Value = 4433
$ws.Range("M94").Value = -2660.25
$ws.Range("N94").Value = -5335
$ws.Range("H98").Value = 1684.7059
$ws.Range("I98").Value = 1402.8572
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 1402.8572
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 95.14280000000008
$ws.Range("N98").Value = -5996
$ws.Range("H101").Value = 8456.15
$ws.Range("J101").Value = 10423.4375
$ws.Range("L101").Value = 31270.3125
$ws.Range("N101").Value = -34514.3125
$ws.Range("H113").Value = 1999.3334
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1999.3334
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1999.3334
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8507.3334
$ws.Range("H122").Value = 1684.7059
$ws.Range("I122").Value = 1402.8572
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4208.571599999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1758.571599999999
$ws.Range("N122").Value = -13900
$ws.Range("H137").Value = 33944
$ws.Range("I137").Value = 1605.3
$ws.Range("J137").Value = 87841.836
$ws.Range("K137").Value = 4815.9
$ws.Range("L137").Value = 263525.508
$ws.Range("M137").Value = -2265.9
$ws.Range("N137").Value = -268625.508

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 50000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 50000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 50000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -50288
$ws.Range("H112").Value = 19375
$ws.Range("J112").Value = 19375
$ws.Range("L112").Value = 19375
$ws.Range("N112").Value = -22329
$ws.Range("H125").Value = 71143.336
$ws.Range("J125").Value = 71143.336
$ws.Range("L125").Value = 71143.336
$ws.Range("N125").Value = -80983.336
$ws.Range("H132").Value = 1121.5319
$ws.Range("I132").Value = 974.7442
$ws.Range("J132").Value = 2699.5
$ws.Range("K132").Value = 2924.2326
$ws.Range("L132").Value = 8098.5
$ws.Range("M132").Value = -394.2325999999998
$ws.Range("N132").Value = -13158.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 642.3333
$ws.Range("I11").Value = 428
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 428
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -288
$ws.Range("N11").Value = -3280
$ws.Range("H115").Value = 22945.857
$ws.Range("I115").Value = 50621
$ws.Range("J115").Value = 18333.334
$ws.Range("K115").Value = 50621
$ws.Range("L115").Value = 18333.334
$ws.Range("N115").Value = -21467.334
$ws.Range("M115").Value = -49054

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2988
$ws.Range("I31").Value = 1801.3
$ws.Range("J31").Value = 4471.375
$ws.Range("K31").Value = 1801.3
$ws.Range("L31").Value = 4471.375
$ws.Range("M31").Value = -1506.3
$ws.Range("N31").Value = -5061.375
$ws.Range("H34").Value = 2988
$ws.Range("I34").Value = 1801.3
$ws.Range("J34").Value = 4471.375
$ws.Range("K34").Value = 1801.3
$ws.Range("L34").Value = 4471.375
$ws.Range("M34").Value = -1599.3
$ws.Range("N34").Value = -4875.375
$ws.Range("H58").Value = 1202.579
$ws.Range("I58").Value = 714.9318
$ws.Range("J58").Value = 2853.077
$ws.Range("K58").Value = 714.9318
$ws.Range("L58").Value = 2853.077
$ws.Range("M58").Value = -511.9318
$ws.Range("N58").Value = -3259.077
$ws.Range("H105").Value = 1340
$ws.Range("I105").Value = 1313.3334
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1313.3334
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 433.6666
$ws.Range("N105").Value = -4994
$ws.Range("H132").Value = 1332
$ws.Range("I132").Value = 1101.8462
$ws.Range("J132").Value = 2828
$ws.Range("K132").Value = 3305.5386
$ws.Range("L132").Value = 8484
$ws.Range("M132").Value = -775.5385999999999
$ws.Range("N132").Value = -13544
$ws.Range("H134").Value = 1222.2826
$ws.Range("I134").Value = 743.8293
$ws.Range("J134").Value = 5145.6
$ws.Range("K134").Value = 2231.4879
$ws.Range("L134").Value = 15436.8
$ws.Range("M134").Value = 303.5120999999999
$ws.Range("N134").Value = -20506.8
$ws.Range("H136").Value = 1202.579
$ws.Range("I136").Value = 714.9318
$ws.Range("J136").Value = 2853.077
$ws.Range("K136").Value = 2144.7954
$ws.Range("L136").Value = 8559.231
$ws.Range("M136").Value = 405.2046
$ws.Range("N136").Value = -13659.231

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 542
$ws.Range("I18").Value = 135.33333
$ws.Range("J18").Value = 1274
$ws.Range("K18").Value = 405.99999
$ws.Range("L18").Value = 3822
$ws.Range("M18").Value = -236.99999
$ws.Range("N18").Value = -4160
$ws.Range("H32").Value = 55556556
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6566
$ws.Range("H39").Value = 3285.2856
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3285.2856
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 9855.856800000001
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -10443.8568
$ws.Range("H103").Value = 1000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 3000
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -4758
$ws.Range("H117").Value = 4582.5713
$ws.Range("I117").Value = 3015.6
$ws.Range("J117").Value = 8500
$ws.Range("K117").Value = 9046.799999999999
$ws.Range("L117").Value = 25500
$ws.Range("M117").Value = -5604.799999999999
$ws.Range("N117").Value = -32384
$ws.Range("H131").Value = 777.7826
$ws.Range("I131").Value = 508.42856
$ws.Range("J131").Value = 895.625
$ws.Range("K131").Value = 1525.28568
$ws.Range("L131").Value = 2686.875
$ws.Range("M131").Value = 3514.71432
$ws.Range("N131").Value = -12766.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 733.3333
$ws.Range("I9").Value = 733.3333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 733.3333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -563.3333
$ws.Range("N9").ClearContents()
$ws.Range("H97").Value = 698.2
$ws.Range("I97").Value = 795
$ws.Range("J97").Value = 311
$ws.Range("K97").Value = 795
$ws.Range("L97").Value = 311
$ws.Range("M97").Value = -299
$ws.Range("N97").Value = -1303
$ws.Range("H102").Value = 1227.2941
$ws.Range("I102").Value = 1215.625
$ws.Range("K102").Value = 1215.625
$ws.Range("M102").Value = 406.375
$ws.Range("H113").Value = 7652.706
$ws.Range("I113").Value = 1485
$ws.Range("J113").Value = 13135.111
$ws.Range("K113").Value = 1485
$ws.Range("L113").Value = 13135.111
$ws.Range("M113").Value = 685
$ws.Range("N113").Value = -17475.111
$ws.Range("H122").Value = 11112871
$ws.Range("I122").Value = 14287121
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 42861363
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -42858913
$ws.Range("N122").Value = -13885

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -683
$ws.Range("N32").Value = -1634
$ws.Range("H122").Value = 3447.4092
$ws.Range("I122").Value = 3447.4092
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10342.2276
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7892.2276
$ws.Range("N122").ClearContents()
